$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.365.11'
$ws.Range('E2').Value = '  +0.14%  '
$ws.Range('D3').Value = '1.937.23'
$ws.Range('E3').Value = '  +0.11%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9996'
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7659'
$ws.Range('E5').Value = '  +6.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '248.66'
$ws.Range('E6').Value = '  -0.64%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9984'
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '28.24'
$ws.Range('E8').Value = '  +1.88%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3215'
$ws.Range('E9').Value = '  -2.63%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07125'
$ws.Range('E10').Value = '  -0.89%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7881'
$ws.Range('E11').Value = '  -2.54%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08003'
$ws.Range('E12').Value = '  -1.07%  '
$ws.Range('D13').Value = '1.934.70'
$ws.Range('E13').Value = '  -0.04%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.383'
$ws.Range('E14').Value = '  -1.86%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '94.86'
$ws.Range('E15').Value = '  +0.55%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.73'
$ws.Range('E16').Value = '  -2.87%  '
$ws.Range('D17').Value = '30.377.84'
$ws.Range('E17').Value = '  +0.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '253.90'
$ws.Range('E18').Value = '  +1.66%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008036'
$ws.Range('E19').Value = '  -2.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.816'
$ws.Range('E20').Value = '  -1.35%  '
$ws.Range('D21').Value = '2.191.93'
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9983'
$ws.Range('E22').Value = '  -0.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').Value = '  -0.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.824'
$ws.Range('E24').Value = '  -2.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.635'
$ws.Range('E25').Value = '  -1.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.83'
$ws.Range('E26').Value = '  +0.79%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1355'
$ws.Range('E27').Value = '  +2.82%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.16'
$ws.Range('E28').Value = '  -0.33%  '
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.307'
$ws.Range('E29').Value = '  -2.76%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.363'
$ws.Range('E30').Value = '  +0.72%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.530'
$ws.Range('E31').Value = '  -2.28%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.443'
$ws.Range('E32').Value = '  +0.42%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.156'
$ws.Range('E33').Value = '  -0.28%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05204'
$ws.Range('E34').Value = '  +0.33%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.290'
$ws.Range('E35').Value = '  +0.91%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7533'
$ws.Range('E36').Value = '  +0.77%  '
$ws.Range('E37').Value = '  +1.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01979'
$ws.Range('E38').Value = '  -0.01%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.809'
$ws.Range('E39').Value = '  -0.59%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '78.51'
$ws.Range('E40').Value = '  -1.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.480'
$ws.Range('E41').Value = '  +1.71%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4528'
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.995'
$ws.Range('E43').Value = '  -1.26%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9993'
$ws.Range('E44').Value = '  -0.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8384'
$ws.Range('E45').Value = '  -0.96%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '102.29'
$ws.Range('E46').Value = '  +0.64%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.856'
$ws.Range('E47').Value = '  +0.99%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.555'
$ws.Range('E48').Value = '  +1.42%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '988.00'
$ws.Range('E49').Value = '  +12.87%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '37.41'
$ws.Range('E50').Value = '  +1.89%  '
$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4178'
$ws.Range('E51').Value = '  +0.02%  '
